{"js": "// Update the date line and every multiplication answer cell in the table,\n// matching the commit that refreshed this worksheet's generated output.\nconst replacements = [\n  [\"2025-09-26 Friday\", \"2025-09-27 Saturday\"],\n  [\"51\u00d734=1734\", \"40\u00d724=960\"],\n  [\"52\u00d753=2756\", \"77\u00d799=7623\"],\n  [\"75\u00d731=2325\", \"79\u00d746=3634\"],\n  [\"54\u00d778=4212\", \"55\u00d779=4345\"],\n  [\"38\u00d795=3610\", \"17\u00d725=425\"],\n  [\"75\u00d759=4425\", \"81\u00d781=6561\"],\n  [\"98\u00d762=6076\", \"70\u00d714=980\"],\n  [\"91\u00d744=4004\", \"56\u00d720=1120\"],\n  [\"20\u00d766=1320\", \"13\u00d740=520\"],\n  [\"95\u00d732=3040\", \"13\u00d769=897\"],\n  [\"18\u00d723=414\", \"74\u00d780=5920\"],\n  [\"46\u00d755=2530\", \"67\u00d732=2144\"],\n  [\"47\u00d794=4418\", \"90\u00d759=5310\"],\n  [\"84\u00d725=2100\", \"41\u00d749=2009\"],\n  [\"36\u00d750=1800\", \"74\u00d794=6956\"],\n  [\"37\u00d715=555\", \"81\u00d757=4617\"],\n  [\"58\u00d780=4640\", \"28\u00d744=1232\"],\n  [\"39\u00d713=507\", \"94\u00d741=3854\"],\n  [\"67\u00d738=2546\", \"22\u00d761=1342\"],\n  [\"15\u00d793=1395\", \"97\u00d779=7663\"],\n  [\"43\u00d755=2365\", \"98\u00d720=1960\"],\n  [\"35\u00d713=455\", \"53\u00d753=2809\"],\n  [\"72\u00d721=1512\", \"34\u00d790=3060\"],\n  [\"14\u00d762=868\", \"37\u00d742=1554\"],\n  [\"78\u00d750=3900\", \"67\u00d767=4489\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and every multiplication answer cell in the table,\n# matching the commit that refreshed this worksheet's generated output.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{Old = \"2025-09-26 Friday\"; New = \"2025-09-27 Saturday\"},\n    @{Old = \"51\u00d734=1734\"; New = \"40\u00d724=960\"},\n    @{Old = \"52\u00d753=2756\"; New = \"77\u00d799=7623\"},\n    @{Old = \"75\u00d731=2325\"; New = \"79\u00d746=3634\"},\n    @{Old = \"54\u00d778=4212\"; New = \"55\u00d779=4345\"},\n    @{Old = \"38\u00d795=3610\"; New = \"17\u00d725=425\"},\n    @{Old = \"75\u00d759=4425\"; New = \"81\u00d781=6561\"},\n    @{Old = \"98\u00d762=6076\"; New = \"70\u00d714=980\"},\n    @{Old = \"91\u00d744=4004\"; New = \"56\u00d720=1120\"},\n    @{Old = \"20\u00d766=1320\"; New = \"13\u00d740=520\"},\n    @{Old = \"95\u00d732=3040\"; New = \"13\u00d769=897\"},\n    @{Old = \"18\u00d723=414\"; New = \"74\u00d780=5920\"},\n    @{Old = \"46\u00d755=2530\"; New = \"67\u00d732=2144\"},\n    @{Old = \"47\u00d794=4418\"; New = \"90\u00d759=5310\"},\n    @{Old = \"84\u00d725=2100\"; New = \"41\u00d749=2009\"},\n    @{Old = \"36\u00d750=1800\"; New = \"74\u00d794=6956\"},\n    @{Old = \"37\u00d715=555\"; New = \"81\u00d757=4617\"},\n    @{Old = \"58\u00d780=4640\"; New = \"28\u00d744=1232\"},\n    @{Old = \"39\u00d713=507\"; New = \"94\u00d741=3854\"},\n    @{Old = \"67\u00d738=2546\"; New = \"22\u00d761=1342\"},\n    @{Old = \"15\u00d793=1395\"; New = \"97\u00d779=7663\"},\n    @{Old = \"43\u00d755=2365\"; New = \"98\u00d720=1960\"},\n    @{Old = \"35\u00d713=455\"; New = \"53\u00d753=2809\"},\n    @{Old = \"72\u00d721=1512\"; New = \"34\u00d790=3060\"},\n    @{Old = \"14\u00d762=868\"; New = \"37\u00d742=1554\"},\n    @{Old = \"78\u00d750=3900\"; New = \"67\u00d767=4489\"}\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.Old\n    $find.Replacement.Text = $pair.New\n    $find.Execute($pair.Old, $false, $false, $false, $false, $false, $true, 1, $false, $pair.New, 2)\n}\n"}
